$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 340.7143
$ws.Range("I38").Value = 188.18182
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 564.5454599999999
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -192.5454599999999
$ws.Range("N38").Value = -3444

$ws.Range("H61").Value = 905.1429000000001
$ws.Range("I61").Value = 905.1429000000001
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2715.4287
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2543.4287
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H69").Value = 3482
$ws.Range("I69").Value = 3956.5
$ws.Range("J69").Value = 3244.75
$ws.Range("K69").Value = 11869.5
$ws.Range("L69").Value = 9734.25
$ws.Range("M69").Value = -10995.5
$ws.Range("N69").Value = -11482.25

$ws.Range("H70").Value = 1728.2858
$ws.Range("I70").Value = 1679.8
$ws.Range("J70").Value = 1849.5
$ws.Range("K70").Value = 5039.4
$ws.Range("L70").Value = 5548.5
$ws.Range("M70").Value = -4769.4
$ws.Range("N70").Value = -6088.5

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H72").Value = 3482
$ws.Range("I72").Value = 3956.5
$ws.Range("J72").Value = 3244.75
$ws.Range("K72").Value = 35608.5
$ws.Range("L72").Value = 29202.75
$ws.Range("M72").Value = -31240.5
$ws.Range("N72").Value = -37938.75

$ws.Range("H73").Value = 1728.2858
$ws.Range("I73").Value = 1679.8
$ws.Range("J73").Value = 1849.5
$ws.Range("K73").Value = 5039.4
$ws.Range("L73").Value = 5548.5
$ws.Range("M73").Value = -4103.4
$ws.Range("N73").Value = -7420.5

$ws.Range("H92").Value = 1938.3889
$ws.Range("I92").Value = 1592.7333
$ws.Range("J92").Value = 3666.6667
$ws.Range("K92").Value = 1592.7333
$ws.Range("L92").Value = 3666.6667
$ws.Range("M92").Value = -344.7333000000001
$ws.Range("N92").Value = -6162.6667

$ws.Range("H96").Value = 1029
$ws.Range("I96").Value = 514
$ws.Range("J96").Value = 1715.6666
$ws.Range("K96").Value = 1542
$ws.Range("L96").Value = 5146.9998
$ws.Range("M96").Value = -169
$ws.Range("N96").Value = -7892.9998

$ws.Range("H106").Value = 100001980
$ws.Range("I106").Value = 120001580
$ws.Range("K106").Value = 120001580
$ws.Range("M106").Value = -120000949

$ws.Range("H138").Value = 2148.8823
$ws.Range("I138").Value = 2270
$ws.Range("J138").Value = 2101.9185
$ws.Range("K138").Value = 6810
$ws.Range("L138").Value = 6305.755500000001
$ws.Range("M138").Value = -1670
$ws.Range("N138").Value = -16585.7555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4791.2
$ws.Range("I45").Value = 4728
$ws.Range("J45").Value = 4833.3335
$ws.Range("K45").Value = 4728
$ws.Range("L45").Value = 4833.3335
$ws.Range("M45").Value = -4351
$ws.Range("N45").Value = -5587.3335

$ws.Range("H122").Value = 58469.445
$ws.Range("I122").Value = 102222.4
$ws.Range("J122").Value = 3778.25
$ws.Range("K122").Value = 306667.2
$ws.Range("L122").Value = 11334.75
$ws.Range("M122").Value = -304217.2
$ws.Range("N122").Value = -16234.75

$ws.Range("H132").Value = 2111.6287
$ws.Range("I132").Value = 1754.037
$ws.Range("J132").Value = 3318.5
$ws.Range("K132").Value = 5262.111
$ws.Range("L132").Value = 9955.5
$ws.Range("M132").Value = -2732.111
$ws.Range("N132").Value = -15015.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6021.6733
$ws.Range("I31").Value = 1430.2963
$ws.Range("J31").Value = 11656.546
$ws.Range("K31").Value = 1430.2963
$ws.Range("L31").Value = 11656.546
$ws.Range("M31").Value = -1135.2963
$ws.Range("N31").Value = -12246.546

$ws.Range("H34").Value = 6021.6733
$ws.Range("I34").Value = 1430.2963
$ws.Range("J34").Value = 11656.546
$ws.Range("K34").Value = 1430.2963
$ws.Range("L34").Value = 11656.546
$ws.Range("M34").Value = -1228.2963
$ws.Range("N34").Value = -12060.546

$ws.Range("H132").Value = 6412185.5
$ws.Range("I132").Value = 1559
$ws.Range("J132").Value = 13891249
$ws.Range("K132").Value = 4677
$ws.Range("L132").Value = 41673747
$ws.Range("M132").Value = -2147
$ws.Range("N132").Value = -41678807

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1015
$ws.Range("I5").Value = 1030
$ws.Range("K5").Value = 3090
$ws.Range("M5").Value = -2978

$ws.Range("H62").Value = 7666.5557
$ws.Range("J62").Value = 7666.5557
$ws.Range("L62").Value = 22999.6671
$ws.Range("N62").Value = -24371.6671

$ws.Range("H63").Value = 4520
$ws.Range("I63").Value = 3800
$ws.Range("J63").Value = 4700
$ws.Range("K63").Value = 11400
$ws.Range("L63").Value = 14100
$ws.Range("M63").Value = -10651
$ws.Range("N63").Value = -15598

$ws.Range("H65").Value = 7666.5557
$ws.Range("J65").Value = 7666.5557
$ws.Range("L65").Value = 68999.0013
$ws.Range("N65").Value = -75863.0013

$ws.Range("H66").Value = 4520
$ws.Range("I66").Value = 3800
$ws.Range("J66").Value = 4700
$ws.Range("K66").Value = 34200
$ws.Range("L66").Value = 42300
$ws.Range("M66").Value = -30456
$ws.Range("N66").Value = -49788

$ws.Range("H70").Value = 1991.8235
$ws.Range("I70").Value = 984.55554
$ws.Range("J70").Value = 3125
$ws.Range("K70").Value = 2953.66662
$ws.Range("L70").Value = 9375
$ws.Range("M70").Value = -2638.66662
$ws.Range("N70").Value = -10005

$ws.Range("H73").Value = 1991.8235
$ws.Range("I73").Value = 984.55554
$ws.Range("J73").Value = 3125
$ws.Range("K73").Value = 2953.66662
$ws.Range("L73").Value = 9375
$ws.Range("M73").Value = -1861.66662
$ws.Range("N73").Value = -11559

$ws.Range("H107").Value = 685.86664
$ws.Range("I107").Value = 672.5
$ws.Range("K107").Value = 2017.5
$ws.Range("M107").Value = -97.5

$ws.Range("H122").Value = 6767.0625
$ws.Range("I122").Value = 359.6154
$ws.Range("K122").Value = 3236.5386
$ws.Range("M122").Value = -786.5386000000003

$ws.Range("H135").Value = 1015
$ws.Range("I135").Value = 1030
$ws.Range("K135").Value = 9270
$ws.Range("M135").Value = -6735

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 11130.5
$ws.Range("J39").Value = 11130.5
$ws.Range("L39").Value = 11130.5
$ws.Range("N39").Value = -12194.5

$ws.Range("H122").Value = 4101.8667
$ws.Range("I122").Value = 3740.125
$ws.Range("J122").Value = 4515.2856
$ws.Range("K122").Value = 11220.375
$ws.Range("L122").Value = 13545.8568
$ws.Range("M122").Value = -8770.375
$ws.Range("N122").Value = -18445.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 44005.25
$ws.Range("J21").Value = 44005.25
$ws.Range("L21").Value = 44005.25
$ws.Range("N21").Value = -44353.25

$ws.Range("H98").Value = 61566.668
$ws.Range("J98").Value = 61566.668
$ws.Range("L98").Value = 61566.668
$ws.Range("N98").Value = -67556.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 18334.5
$ws.Range("J15").Value = 18334.5
$ws.Range("L15").Value = 18334.5
$ws.Range("N15").Value = -18910.5

$ws.Range("H19").Value = 32504.5
$ws.Range("J19").Value = 32504.5
$ws.Range("L19").Value = 32504.5
$ws.Range("N19").Value = -32852.5

$ws.Range("H132").Value = 3705905
$ws.Range("I132").Value = 2553.353
$ws.Range("J132").Value = 5954368.5
$ws.Range("K132").Value = 7660.059
$ws.Range("L132").Value = 17863105.5
$ws.Range("M132").Value = -5130.059
$ws.Range("N132").Value = -17868165.5
